$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 3 rows for MuSCs-as-sender (rows 8-10); data was recomputed without that source cluster
$ws.Range("A8:T10").EntireRow.Delete() | Out-Null

# Refresh all remaining data rows (2-7) with the recalculated TPM-derived metrics
# Row 2: ECs -> ECs (Tnf/Tnfrsf1a)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tnf"
$ws.Cells.Item(2, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 2.913576333333333
$ws.Cells.Item(2, 8).Value = 8.740729
$ws.Cells.Item(2, 9).Value = 0.8649322955011439
$ws.Cells.Item(2, 10).Value = 0.8649322955011439
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 14.65767833333333
$ws.Cells.Item(2, 14).Value = 43.973035
$ws.Cells.Item(2, 15).Value = 0.2345581433878666
$ws.Cells.Item(2, 16).Value = 0.2345581433878665
$ws.Cells.Item(2, 17).Value = 42.70626469361278
$ws.Cells.Item(2, 18).Value = 384.356382242515
$ws.Cells.Item(2, 19).Value = 0.2028769133889539
$ws.Cells.Item(2, 20).Value = 0.2028769133889539

# Row 3: ECs -> FAPs (Tnf/Tnfrsf1a)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tnf"
$ws.Cells.Item(3, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 2.913576333333333
$ws.Cells.Item(3, 8).Value = 8.740729
$ws.Cells.Item(3, 9).Value = 0.8649322955011439
$ws.Cells.Item(3, 10).Value = 0.8649322955011439
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 31.695371
$ws.Cells.Item(3, 14).Value = 95.086113
$ws.Cells.Item(3, 15).Value = 0.5072022462686253
$ws.Cells.Item(3, 16).Value = 0.5072022462686253
$ws.Cells.Item(3, 17).Value = 92.34688282181966
$ws.Cells.Item(3, 18).Value = 831.121945396377
$ws.Cells.Item(3, 19).Value = 0.4386956031484586
$ws.Cells.Item(3, 20).Value = 0.4386956031484586

# Row 4: ECs -> MuSCs (Tnf/Tnfrsf1a)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tnf"
$ws.Cells.Item(4, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2.913576333333333
$ws.Cells.Item(4, 8).Value = 8.740729
$ws.Cells.Item(4, 9).Value = 0.8649322955011439
$ws.Cells.Item(4, 10).Value = 0.8649322955011439
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 16.13754733333333
$ws.Cells.Item(4, 14).Value = 48.41264200000001
$ws.Cells.Item(4, 15).Value = 0.2582396103435082
$ws.Cells.Item(4, 16).Value = 0.2582396103435082
$ws.Cells.Item(4, 17).Value = 47.01797598844644
$ws.Cells.Item(4, 18).Value = 423.1617838960181
$ws.Cells.Item(4, 19).Value = 0.2233597789637315
$ws.Cells.Item(4, 20).Value = 0.2233597789637315

# Row 5: FAPs -> ECs (Tnf/Tnfrsf1a)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Tnf"
$ws.Cells.Item(5, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 0.4549836666666667
$ws.Cells.Item(5, 8).Value = 1.364951
$ws.Cells.Item(5, 9).Value = 0.1350677044988561
$ws.Cells.Item(5, 10).Value = 0.1350677044988561
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 14.65767833333333
$ws.Cells.Item(5, 14).Value = 43.973035
$ws.Cells.Item(5, 15).Value = 0.2345581433878666
$ws.Cells.Item(5, 16).Value = 0.2345581433878665
$ws.Cells.Item(5, 17).Value = 6.669004232920557
$ws.Cells.Item(5, 18).Value = 60.021038096285
$ws.Cells.Item(5, 19).Value = 0.03168122999891267
$ws.Cells.Item(5, 20).Value = 0.03168122999891268

# Row 6: FAPs -> FAPs (Tnf/Tnfrsf1a)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tnf"
$ws.Cells.Item(6, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 2
$ws.Cells.Item(6, 6).Value = 0.6666666666666666
$ws.Cells.Item(6, 7).Value = 0.4549836666666667
$ws.Cells.Item(6, 8).Value = 1.364951
$ws.Cells.Item(6, 9).Value = 0.1350677044988561
$ws.Cells.Item(6, 10).Value = 0.1350677044988561
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 31.695371
$ws.Cells.Item(6, 14).Value = 95.086113
$ws.Cells.Item(6, 15).Value = 0.5072022462686253
$ws.Cells.Item(6, 16).Value = 0.5072022462686253
$ws.Cells.Item(6, 17).Value = 14.42087611394033
$ws.Cells.Item(6, 18).Value = 129.787885025463
$ws.Cells.Item(6, 19).Value = 0.06850664312016672
$ws.Cells.Item(6, 20).Value = 0.06850664312016673

# Row 7: FAPs -> MuSCs (Tnf/Tnfrsf1a)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tnf"
$ws.Cells.Item(7, 3).Value = "Tnfrsf1a"
$ws.Cells.Item(7, 4).Value = "MuSCs"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 0.6666666666666666
$ws.Cells.Item(7, 7).Value = 0.4549836666666667
$ws.Cells.Item(7, 8).Value = 1.364951
$ws.Cells.Item(7, 9).Value = 0.1350677044988561
$ws.Cells.Item(7, 10).Value = 0.1350677044988561
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 16.13754733333333
$ws.Cells.Item(7, 14).Value = 48.41264200000001
$ws.Cells.Item(7, 15).Value = 0.2582396103435082
$ws.Cells.Item(7, 16).Value = 0.2582396103435082
$ws.Cells.Item(7, 17).Value = 7.342320456726889
$ws.Cells.Item(7, 18).Value = 66.080884110542
$ws.Cells.Item(7, 19).Value = 0.03487983137977671
$ws.Cells.Item(7, 20).Value = 0.03487983137977672
